# Insert a new data row before the existing row 273 on the "Zapallo" sheet.
# This pushes the previous rows 273-351 down to become rows 274-352 and
# grows the used range from A1:R351 to A1:R352.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(273).Insert()

# Populate the newly inserted row 273 with the new record.
$ws.Range("A273").Value = 10
$ws.Range("B273").Value = "Vega Modelo de Temuco"
$ws.Range("C273").Value = "La Araucanía"
$ws.Range("D273").Value = 44463
$ws.Range("E273").Value = 9
$ws.Range("F273").Value = 100112045
$ws.Range("G273").Value = "Zapallo"
$ws.Range("H273").Value = "Paine"
$ws.Range("I273").Value = "1a (guarda)"
$ws.Range("J273").Value = 1300
$ws.Range("K273").Value = 500
$ws.Range("L273").Value = 600
$ws.Range("M273").Value = 554
$ws.Range("N273").Value = "$/kilo (volumen en unidades)"
$ws.Range("O273").Value = "Región del Maule"
$ws.Range("P273").Value = 554
$ws.Range("Q273").Value = 1
$ws.Range("R273").Value = "Hortaliza"
